# Updates cryptos list cell values per the diff (row-by-row re-ranking refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Use a leading apostrophe so Excel stores numeric-looking strings as text,
    # then restore the Normal style so no stray quote-prefix formatting sticks.
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "71.342.36"
Set-TextValue $ws.Range("E2") "  +0.52%  "

Set-TextValue $ws.Range("D3") "3.842.01"
Set-TextValue $ws.Range("E3") "  +1.07%  "

Set-TextValue $ws.Range("E4") "  +0.03%  "

Set-TextValue $ws.Range("D5") "715.39"
Set-TextValue $ws.Range("E5") "  +1.77%  "

Set-TextValue $ws.Range("D6") "173.12"
Set-TextValue $ws.Range("E6") "  +0.03%  "

Set-TextValue $ws.Range("D7") "3.843.01"
Set-TextValue $ws.Range("E7") "  +1.11%  "

Set-TextValue $ws.Range("E8") "  -0.02%  "

Set-TextValue $ws.Range("E9") "  -0.07%  "

Set-TextValue $ws.Range("E10") "  +0.44%  "

Set-TextValue $ws.Range("D11") "7.37"
Set-TextValue $ws.Range("E11") "  +1.52%  "

Set-TextValue $ws.Range("D12") "0.461"
Set-TextValue $ws.Range("E12") "  +0.14%  "

Set-TextValue $ws.Range("E13") "  -0.20%  "

Set-TextValue $ws.Range("D14") "36.89"
Set-TextValue $ws.Range("E14") "  +2.23%  "

Set-TextValue $ws.Range("D15") "4.488.08"
Set-TextValue $ws.Range("E15") "  +1.05%  "

Set-TextValue $ws.Range("D16") "3.833.23"
Set-TextValue $ws.Range("E16") "  +1.13%  "

Set-TextValue $ws.Range("D17") "71.274.32"
Set-TextValue $ws.Range("E17") "  +0.56%  "

Set-TextValue $ws.Range("D18") "7.26"
Set-TextValue $ws.Range("E18") "  +0.77%  "

Set-TextValue $ws.Range("E19") "  +0.45%  "

Set-TextValue $ws.Range("D20") "17.45"
Set-TextValue $ws.Range("E20") "  -1.16%  "

Set-TextValue $ws.Range("B21") "BitcoinCash"
Set-TextValue $ws.Range("C21") "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue $ws.Range("D21") "497.69"
Set-TextValue $ws.Range("E21") "  +3.44%  "

Set-TextValue $ws.Range("B22") "Uniswap"
Set-TextValue $ws.Range("C22") "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue $ws.Range("D22") "10.73"
Set-TextValue $ws.Range("E22") "  -2.05%  "

Set-TextValue $ws.Range("D23") "0.733"
Set-TextValue $ws.Range("E23") "  +2.81%  "

Set-TextValue $ws.Range("D24") "85.32"
Set-TextValue $ws.Range("E24") "  +1.32%  "

Set-TextValue $ws.Range("E25") "  +1.91%  "

Set-TextValue $ws.Range("D26") "10.70"
Set-TextValue $ws.Range("E26") "  +1.80%  "

Set-TextValue $ws.Range("D27") "12.20"
Set-TextValue $ws.Range("E27") "  -0.81%  "

Set-TextValue $ws.Range("E28") "  -2.56%  "

Set-TextValue $ws.Range("D29") "3.18"
Set-TextValue $ws.Range("E29") "  +1.61%  "

Set-TextValue $ws.Range("E30") "  -0.02%  "

Set-TextValue $ws.Range("E31") "  -0.84%  "

Set-TextValue $ws.Range("E32") "  -2.12%  "

Set-TextValue $ws.Range("D33") "29.48"

Set-TextValue $ws.Range("D34") "0.181"
Set-TextValue $ws.Range("E34") "  -4.50%  "

Set-TextValue $ws.Range("D35") "9.23"
Set-TextValue $ws.Range("E35") "  -0.42%  "

Set-TextValue $ws.Range("D36") "3.805.45"
Set-TextValue $ws.Range("E36") "  +1.44%  "

Set-TextValue $ws.Range("D37") "0.998"
Set-TextValue $ws.Range("E37") "  -0.21%  "

Set-TextValue $ws.Range("E38") "  +0.56%  "

Set-TextValue $ws.Range("E39") "  +5.55%  "

Set-TextValue $ws.Range("D40") "6.03"
Set-TextValue $ws.Range("E40") "  +0.40%  "

Set-TextValue $ws.Range("E41") "  -1.63%  "

Set-TextValue $ws.Range("E42") "  +2.69%  "

Set-TextValue $ws.Range("E44") "  +0.09%  "

Set-TextValue $ws.Range("D45") "0.000322"
Set-TextValue $ws.Range("E45") "  +0.39%  "

Set-TextValue $ws.Range("D46") "163.76"
Set-TextValue $ws.Range("E46") "  -0.22%  "

Set-TextValue $ws.Range("D47") "48.96"
Set-TextValue $ws.Range("E47") "  +0.17%  "

Set-TextValue $ws.Range("D48") "426.12"
Set-TextValue $ws.Range("E48") "  +3.55%  "

Set-TextValue $ws.Range("B49") "Cosmos"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D49") "8.74"
Set-TextValue $ws.Range("E49") "  +1.24%  "

Set-TextValue $ws.Range("B50") "ONDO"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-TextValue $ws.Range("D50") "1.39"
Set-TextValue $ws.Range("E50") "  -0.09%  "

Set-TextValue $ws.Range("E51") "  -1.04%  "
